# PlayerPerformance_4428.xlsx — additional scraping edit
#
# NOTE on this COM shim: a worksheet variable captured via Worksheets.Item(n)
# tracks the *index*, not the object identity, so after any Worksheets.Add()
# call shifts indices around, previously-captured sheet variables can end up
# pointing at the wrong sheet. To stay safe we always re-fetch the sheet we
# want to touch by *name* immediately before using it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert "Player Info" immediately before "ODI Batting".
# ---------------------------------------------------------------------------
$battingRef = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingRef)
$playerInfo.Name = "Player Info"

# ---------------------------------------------------------------------------
# 2) Insert "ODI Batting Extra" immediately after "ODI Batting".
# ---------------------------------------------------------------------------
$battingRef2 = $wb.Worksheets.Item("ODI Batting")
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $battingRef2)
$battingExtra.Name = "ODI Batting Extra"

# ---------------------------------------------------------------------------
# 3) Populate "Player Info".
# ---------------------------------------------------------------------------
$playerInfoSheet = $wb.Worksheets.Item("Player Info")

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $playerInfoHeaders.Length; $i++) {
    $cell = $playerInfoSheet.Cells.Item(1, $i + 1)
    $cell.Value = $playerInfoHeaders[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$playerInfoRow = @("4428", "Samuel William Billings", "Right Handed", "Does Not Bowl | Unknown")
for ($i = 0; $i -lt $playerInfoRow.Length; $i++) {
    $playerInfoSheet.Cells.Item(2, $i + 1).Value = $playerInfoRow[$i]
}

# ---------------------------------------------------------------------------
# 4) Populate "ODI Batting Extra".
# ---------------------------------------------------------------------------
$battingExtraSheet = $wb.Worksheets.Item("ODI Batting Extra")

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $extraHeaders.Length; $i++) {
    $cell = $battingExtraSheet.Cells.Item(1, $i + 1)
    $cell.Value = $extraHeaders[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$extraRows = @(
    @("4004", 2, "0", "0", "", "NO"),
    @("4020", 6, "", "", "", "NO"),
    @("4021", 6, "1", "0", "2.13%", "NO"),
    @("4073", 6, "0", "0", "0.78%", "NO"),
    @("4075", 6, "", "", "", "NO"),
    @("4165", 6, "1", "0", "3.29%", "NO"),
    @("4167", "", "", "", "", "NO"),
    @("4426", "", "", "", "", "NO"),
    @("4427", 5, "6", "0", "21.30%", "NO"),
    @("4428", 6, "3", "0", "5.79%", "NO"),
    @("4429", 6, "14", "2", "42.91%", "NO"),
    @("4430", 6, "0", "0", "3.46%", "NO"),
    @("4431", 6, "4", "2", "18.87%", "NO"),
    @("4454", 6, "1", "0", "7.17%", "NO"),
    @("4469", 5, "0", "0", "1.59%", "NO"),
    @("4470", 5, "", "", "", "NO"),
    @("4471", "", "", "", "", "NO"),
    @("4660", "", "", "", "", "NO"),
    @("4663", "", "", "", "", "NO"),
    @("4666", "", "", "", "", "NO")
)

for ($r = 0; $r -lt $extraRows.Length; $r++) {
    $row = $extraRows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $battingExtraSheet.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------------
# 5) Update "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, replace the
#    URL values with the bare match code, and clear the already-empty
#    INNING_NUMBER cells on the "did not bat" rows.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$matchCodes = @("3804","3805","3806","3807","3809","3949","3978","4001","4004","4020","4021","4073","4075","4165","4167","4426","4427","4428","4429","4430","4431","4454","4469","4470","4471","4660","4663","4666")
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $battingSheet.Cells.Item($i + 2, 4).Value = $matchCodes[$i]
}

$blankInningRows = @(5, 11, 14, 25, 26)
foreach ($r in $blankInningRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}
